$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster

$smTcs = $sm.Theme.ThemeColorScheme
$nmTcs = $nm.Theme.ThemeColorScheme

Write-Host "SM c3 before: $($smTcs.Item(3).RGB)"
Write-Host "NM c3 before: $($nmTcs.Item(3).RGB)"

$smTcs.Item(3).RGB = 111111
Write-Host "SM c3 after set: $($smTcs.Item(3).RGB)"
Write-Host "NM c3 after SM set (should be unaffected if separate): $($nmTcs.Item(3).RGB)"

$nmTcs.Item(4).RGB = 222222
Write-Host "NM c4 after set: $($nmTcs.Item(4).RGB)"
Write-Host "SM c4 after NM set (should be unaffected if separate): $($smTcs.Item(4).RGB)"
